$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedules")

# Insert a new row at position 2, pushing the existing rows (old rows 2-30,
# now 3-31) down by one.
$ws.Rows.Item(2).Insert()

# Populate the freshly inserted row with the new "default" schedule entry.
$ws.Range("A2").Value = "default"
$ws.Range("B2").Value = 7
$ws.Range("C2").Value = 19
$ws.Range("D2").Value = 19
$ws.Range("E2").Value = 7

# Reflect the author's final cell selection on the sheet.
$ws.Range("H13").Select()
